$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the existing legend table (E1:F3) to its new location (J1:K3). ---
# This keeps the right-aligned style (style index 1) that was already used
# on the label column, now landing on column J.
$ws.Range("E1:F3").Cut($ws.Range("J1:K3"))

# Remove the leftover (now unused) cells in column E rows 2-3.
$ws.Range("E2:E3").Clear()

# --- New "Client" / "Prestataire" header cells ---
$ws.Range("D1").Value = "Client"

$ws.Range("E1").Value = "Prestataire"
# -4131 = xlLeft (creates/uses the new left-aligned style)
$ws.Range("E1").HorizontalAlignment = -4131

# --- Extend the legend table with the new rows, reusing J1's (right aligned) style ---
$ws.Range("J1").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("J4").Value = "R :"

$ws.Range("J1").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("J5").Value = "I :"

$ws.Range("K4").Value = "Seconde main"
$ws.Range("K5").Value = "Aller une troisième si elle veut au moins elle est au courant"

$ws.Range("J1").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("J6").Value = "C : "
$ws.Range("K6").Value = "Consulter"

# --- Update the legend for "A :" with the new wording ---
$ws.Range("J2").Value = " Au moins un A :"

# --- Column width tweaks ---
# (target stored widths are 14.7109375 / 16.7109375; the runtime quantizes
# ColumnWidth to roughly 1/6-character steps, so we pick the input that
# lands on the closest achievable stored width)
$ws.Columns("F").ColumnWidth = 13.8
$ws.Columns("H").ColumnWidth = 15.8

# --- Update the active selection ---
[void]$ws.Range("I9").Select()
